# Swap the contents of rows 8 and 9 (for the columns that actually differ)
# in the "Artfynd" worksheet. Columns C, I, K, T, U, V, W, Y, AA, AD, AE, AG,
# AT, AW, AX, AY are identical between the two rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R", "S", "Z", "AB")

foreach ($col in $cols) {
    $rng8 = $ws.Range("$col`8")
    $rng9 = $ws.Range("$col`9")

    $val8 = $rng8.Value2
    $val9 = $rng9.Value2

    $rng8.Value2 = $val9
    $rng9.Value2 = $val8
}
